$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$arr = New-Object 'object[,]' 96,14
$arr[0,0] = 45740.0
$arr[0,1] = 0.0
$arr[0,2] = 504.325
$arr[0,3] = 82.0
$arr[0,4] = 0.0
$arr[0,5] = 1163.3
$arr[0,6] = 0.0
$arr[0,7] = 0.0
$arr[0,8] = 700.0
$arr[0,9] = 0.0
$arr[0,10] = 135.0
$arr[0,11] = 1245.3
$arr[0,12] = 1339.325
$arr[0,13] = -94.02500000000009
$arr[1,0] = 45740.01041666666
$arr[1,1] = 0.0
$arr[1,2] = 504.325
$arr[1,3] = 82.0
$arr[1,4] = 0.0
$arr[1,5] = 1297.6
$arr[1,6] = 0.0
$arr[1,7] = 0.0
$arr[1,8] = 700.0
$arr[1,9] = 0.0
$arr[1,10] = 135.0
$arr[1,11] = 1379.6
$arr[1,12] = 1339.325
$arr[1,13] = 40.27499999999986
$arr[2,0] = 45740.02083333334
$arr[2,1] = 0.0
$arr[2,2] = 504.325
$arr[2,3] = 82.0
$arr[2,4] = 0.0
$arr[2,5] = 1250.9
$arr[2,6] = 0.0
$arr[2,7] = 0.0
$arr[2,8] = 700.0
$arr[2,9] = 0.0
$arr[2,10] = 135.0
$arr[2,11] = 1332.9
$arr[2,12] = 1339.325
$arr[2,13] = -6.425000000000182
$arr[3,0] = 45740.03125
$arr[3,1] = 0.0
$arr[3,2] = 504.325
$arr[3,3] = 82.0
$arr[3,4] = 0.0
$arr[3,5] = 1200.0
$arr[3,6] = 0.0
$arr[3,7] = 0.0
$arr[3,8] = 700.0
$arr[3,9] = 0.0
$arr[3,10] = 135.0
$arr[3,11] = 1282.0
$arr[3,12] = 1339.325
$arr[3,13] = -57.32500000000005
$arr[4,0] = 45740.04166666666
$arr[4,1] = 0.0
$arr[4,2] = 534.25
$arr[4,3] = 41.0
$arr[4,4] = 0.0
$arr[4,5] = 1248.6
$arr[4,6] = 0.0
$arr[4,7] = 0.0
$arr[4,8] = 624.0
$arr[4,9] = 0.0
$arr[4,10] = 159.0
$arr[4,11] = 1289.6
$arr[4,12] = 1317.25
$arr[4,13] = -27.65000000000009
$arr[5,0] = 45740.05208333334
$arr[5,1] = 0.0
$arr[5,2] = 534.25
$arr[5,3] = 41.0
$arr[5,4] = 0.0
$arr[5,5] = 1327.0
$arr[5,6] = 0.0
$arr[5,7] = 0.0
$arr[5,8] = 624.0
$arr[5,9] = 0.0
$arr[5,10] = 159.0
$arr[5,11] = 1368.0
$arr[5,12] = 1317.25
$arr[5,13] = 50.75
$arr[6,0] = 45740.0625
$arr[6,1] = 0.0
$arr[6,2] = 534.25
$arr[6,3] = 41.0
$arr[6,4] = 0.0
$arr[6,5] = 1337.0
$arr[6,6] = 0.0
$arr[6,7] = 0.0
$arr[6,8] = 624.0
$arr[6,9] = 0.0
$arr[6,10] = 159.0
$arr[6,11] = 1378.0
$arr[6,12] = 1317.25
$arr[6,13] = 60.75
$arr[7,0] = 45740.07291666666
$arr[7,1] = 0.0
$arr[7,2] = 534.25
$arr[7,3] = 41.0
$arr[7,4] = 0.0
$arr[7,5] = 1316.0
$arr[7,6] = 0.0
$arr[7,7] = 0.0
$arr[7,8] = 624.0
$arr[7,9] = 0.0
$arr[7,10] = 159.0
$arr[7,11] = 1357.0
$arr[7,12] = 1317.25
$arr[7,13] = 39.75
$arr[8,0] = 45740.08333333334
$arr[8,1] = 0.0
$arr[8,2] = 485.625
$arr[8,3] = 19.0
$arr[8,4] = 0.0
$arr[8,5] = 1255.0
$arr[8,6] = 0.0
$arr[8,7] = 0.0
$arr[8,8] = 543.0
$arr[8,9] = 0.0
$arr[8,10] = 143.0
$arr[8,11] = 1274.0
$arr[8,12] = 1171.625
$arr[8,13] = 102.375
$arr[9,0] = 45740.09375
$arr[9,1] = 0.0
$arr[9,2] = 485.625
$arr[9,3] = 19.0
$arr[9,4] = 0.0
$arr[9,5] = 1289.0
$arr[9,6] = 0.0
$arr[9,7] = 0.0
$arr[9,8] = 543.0
$arr[9,9] = 0.0
$arr[9,10] = 143.0
$arr[9,11] = 1308.0
$arr[9,12] = 1171.625
$arr[9,13] = 136.375
$arr[10,0] = 45740.10416666666
$arr[10,1] = 0.0
$arr[10,2] = 485.625
$arr[10,3] = 19.0
$arr[10,4] = 0.0
$arr[10,5] = 1303.6
$arr[10,6] = 0.0
$arr[10,7] = 0.0
$arr[10,8] = 543.0
$arr[10,9] = 0.0
$arr[10,10] = 143.0
$arr[10,11] = 1322.6
$arr[10,12] = 1171.625
$arr[10,13] = 150.9749999999999
$arr[11,0] = 45740.11458333334
$arr[11,1] = 0.0
$arr[11,2] = 485.625
$arr[11,3] = 19.0
$arr[11,4] = 0.0
$arr[11,5] = 1275.8
$arr[11,6] = 0.0
$arr[11,7] = 0.0
$arr[11,8] = 543.0
$arr[11,9] = 0.0
$arr[11,10] = 143.0
$arr[11,11] = 1294.8
$arr[11,12] = 1171.625
$arr[11,13] = 123.175
$arr[12,0] = 45740.125
$arr[12,1] = 0.0
$arr[12,2] = 325.85
$arr[12,3] = 51.0
$arr[12,4] = 0.0
$arr[12,5] = 958.0
$arr[12,6] = 0.0
$arr[12,7] = 0.0
$arr[12,8] = 493.0
$arr[12,9] = 0.0
$arr[12,10] = 83.0
$arr[12,11] = 1009.0
$arr[12,12] = 901.85
$arr[12,13] = 107.15
$arr[13,0] = 45740.13541666666
$arr[13,1] = 0.0
$arr[13,2] = 325.85
$arr[13,3] = 51.0
$arr[13,4] = 0.0
$arr[13,5] = 932.5
$arr[13,6] = 0.0
$arr[13,7] = 0.0
$arr[13,8] = 493.0
$arr[13,9] = 0.0
$arr[13,10] = 83.0
$arr[13,11] = 983.5
$arr[13,12] = 901.85
$arr[13,13] = 81.64999999999998
$arr[14,0] = 45740.14583333334
$arr[14,1] = 0.0
$arr[14,2] = 325.85
$arr[14,3] = 51.0
$arr[14,4] = 0.0
$arr[14,5] = 886.4000000000001
$arr[14,6] = 0.0
$arr[14,7] = 0.0
$arr[14,8] = 493.0
$arr[14,9] = 0.0
$arr[14,10] = 83.0
$arr[14,11] = 937.4000000000001
$arr[14,12] = 901.85
$arr[14,13] = 35.55000000000007
$arr[15,0] = 45740.15625
$arr[15,1] = 0.0
$arr[15,2] = 325.85
$arr[15,3] = 51.0
$arr[15,4] = 0.0
$arr[15,5] = 857.4000000000001
$arr[15,6] = 0.0
$arr[15,7] = 0.0
$arr[15,8] = 493.0
$arr[15,9] = 0.0
$arr[15,10] = 83.0
$arr[15,11] = 908.4000000000001
$arr[15,12] = 901.85
$arr[15,13] = 6.550000000000068
$arr[16,0] = 45740.16666666666
$arr[16,1] = 0.0
$arr[16,2] = 195.575
$arr[16,3] = 189.0
$arr[16,4] = 0.0
$arr[16,5] = 404.2
$arr[16,6] = 0.0
$arr[16,7] = 0.0
$arr[16,8] = 455.0
$arr[16,9] = 13.0
$arr[16,10] = 0.0
$arr[16,11] = 606.2
$arr[16,12] = 650.575
$arr[16,13] = -44.375
$arr[17,0] = 45740.17708333334
$arr[17,1] = 0.0
$arr[17,2] = 195.575
$arr[17,3] = 189.0
$arr[17,4] = 0.0
$arr[17,5] = 353.2
$arr[17,6] = 0.0
$arr[17,7] = 0.0
$arr[17,8] = 455.0
$arr[17,9] = 13.0
$arr[17,10] = 0.0
$arr[17,11] = 555.2
$arr[17,12] = 650.575
$arr[17,13] = -95.375
$arr[18,0] = 45740.1875
$arr[18,1] = 0.0
$arr[18,2] = 195.575
$arr[18,3] = 189.0
$arr[18,4] = 0.0
$arr[18,5] = 381.8
$arr[18,6] = 0.0
$arr[18,7] = 0.0
$arr[18,8] = 455.0
$arr[18,9] = 13.0
$arr[18,10] = 0.0
$arr[18,11] = 583.8
$arr[18,12] = 650.575
$arr[18,13] = -66.77500000000009
$arr[19,0] = 45740.19791666666
$arr[19,1] = 0.0
$arr[19,2] = 195.575
$arr[19,3] = 189.0
$arr[19,4] = 0.0
$arr[19,5] = 262.2
$arr[19,6] = 0.0
$arr[19,7] = 0.0
$arr[19,8] = 455.0
$arr[19,9] = 13.0
$arr[19,10] = 0.0
$arr[19,11] = 464.2
$arr[19,12] = 650.575
$arr[19,13] = -186.375
$arr[20,0] = 45740.20833333334
$arr[20,1] = 52.57499900000005
$arr[20,2] = 0.0
$arr[20,3] = 321.0
$arr[20,4] = 0.0
$arr[20,5] = 0.0
$arr[20,6] = 308.2
$arr[20,7] = 0.0
$arr[20,8] = 260.0
$arr[20,9] = 210.0
$arr[20,10] = 0.0
$arr[20,11] = 583.574999
$arr[20,12] = 568.2
$arr[20,13] = 15.374999
$arr[21,0] = 45740.21875
$arr[21,1] = 52.57499900000005
$arr[21,2] = 0.0
$arr[21,3] = 321.0
$arr[21,4] = 0.0
$arr[21,5] = 0.0
$arr[21,6] = 333.2
$arr[21,7] = 0.0
$arr[21,8] = 260.0
$arr[21,9] = 210.0
$arr[21,10] = 0.0
$arr[21,11] = 583.574999
$arr[21,12] = 593.2
$arr[21,13] = -9.625000999999997
$arr[22,0] = 45740.22916666666
$arr[22,1] = 52.57499900000005
$arr[22,2] = 0.0
$arr[22,3] = 321.0
$arr[22,4] = 0.0
$arr[22,5] = 0.0
$arr[22,6] = 473.1
$arr[22,7] = 0.0
$arr[22,8] = 260.0
$arr[22,9] = 210.0
$arr[22,10] = 0.0
$arr[22,11] = 583.574999
$arr[22,12] = 733.1
$arr[22,13] = -149.525001
$arr[23,0] = 45740.23958333334
$arr[23,1] = 52.57499900000005
$arr[23,2] = 0.0
$arr[23,3] = 321.0
$arr[23,4] = 0.0
$arr[23,5] = 0.0
$arr[23,6] = 546.0
$arr[23,7] = 0.0
$arr[23,8] = 260.0
$arr[23,9] = 210.0
$arr[23,10] = 0.0
$arr[23,11] = 583.574999
$arr[23,12] = 806.0
$arr[23,13] = -222.425001
$arr[24,0] = 45740.25
$arr[24,1] = 0.0
$arr[24,2] = 209.125
$arr[24,3] = 479.0
$arr[24,4] = 0.0
$arr[24,5] = 65.19999999999999
$arr[24,6] = 0.0
$arr[24,7] = 0.0
$arr[24,8] = 342.0
$arr[24,9] = 0.0
$arr[24,10] = 17.0
$arr[24,11] = 544.2
$arr[24,12] = 568.125
$arr[24,13] = -23.92499999999995
$arr[25,0] = 45740.26041666666
$arr[25,1] = 0.0
$arr[25,2] = 209.125
$arr[25,3] = 479.0
$arr[25,4] = 0.0
$arr[25,5] = 95.5
$arr[25,6] = 0.0
$arr[25,7] = 0.0
$arr[25,8] = 342.0
$arr[25,9] = 0.0
$arr[25,10] = 17.0
$arr[25,11] = 574.5
$arr[25,12] = 568.125
$arr[25,13] = 6.375
$arr[26,0] = 45740.27083333334
$arr[26,1] = 0.0
$arr[26,2] = 209.125
$arr[26,3] = 479.0
$arr[26,4] = 0.0
$arr[26,5] = 121.5
$arr[26,6] = 0.0
$arr[26,7] = 0.0
$arr[26,8] = 342.0
$arr[26,9] = 0.0
$arr[26,10] = 17.0
$arr[26,11] = 600.5
$arr[26,12] = 568.125
$arr[26,13] = 32.375
$arr[27,0] = 45740.28125
$arr[27,1] = 0.0
$arr[27,2] = 209.125
$arr[27,3] = 479.0
$arr[27,4] = 0.0
$arr[27,5] = 81.5
$arr[27,6] = 0.0
$arr[27,7] = 0.0
$arr[27,8] = 342.0
$arr[27,9] = 0.0
$arr[27,10] = 17.0
$arr[27,11] = 560.5
$arr[27,12] = 568.125
$arr[27,13] = -7.625
$arr[28,0] = 45740.29166666666
$arr[28,1] = 0.0
$arr[28,2] = 338.2
$arr[28,3] = 425.0
$arr[28,4] = 0.0
$arr[28,5] = 254.0
$arr[28,6] = 0.0
$arr[28,7] = 0.0
$arr[28,8] = 233.0
$arr[28,9] = 0.0
$arr[28,10] = 91.0
$arr[28,11] = 679.0
$arr[28,12] = 662.2
$arr[28,13] = 16.79999999999995
$arr[29,0] = 45740.30208333334
$arr[29,1] = 0.0
$arr[29,2] = 338.2
$arr[29,3] = 425.0
$arr[29,4] = 0.0
$arr[29,5] = 281.3
$arr[29,6] = 0.0
$arr[29,7] = 0.0
$arr[29,8] = 233.0
$arr[29,9] = 0.0
$arr[29,10] = 91.0
$arr[29,11] = 706.3
$arr[29,12] = 662.2
$arr[29,13] = 44.09999999999991
$arr[30,0] = 45740.3125
$arr[30,1] = 0.0
$arr[30,2] = 338.2
$arr[30,3] = 425.0
$arr[30,4] = 0.0
$arr[30,5] = 344.1
$arr[30,6] = 0.0
$arr[30,7] = 0.0
$arr[30,8] = 233.0
$arr[30,9] = 0.0
$arr[30,10] = 91.0
$arr[30,11] = 769.1
$arr[30,12] = 662.2
$arr[30,13] = 106.9
$arr[31,0] = 45740.32291666666
$arr[31,1] = 0.0
$arr[31,2] = 338.2
$arr[31,3] = 425.0
$arr[31,4] = 0.0
$arr[31,5] = 377.9
$arr[31,6] = 0.0
$arr[31,7] = 0.0
$arr[31,8] = 233.0
$arr[31,9] = 0.0
$arr[31,10] = 91.0
$arr[31,11] = 802.9
$arr[31,12] = 662.2
$arr[31,13] = 140.6999999999999
$arr[32,0] = 45740.33333333334
$arr[32,1] = 0.0
$arr[32,2] = 0.0
$arr[32,3] = 0.0
$arr[32,4] = 0.0
$arr[32,5] = 0.0
$arr[32,6] = 0.0
$arr[32,7] = 0.0
$arr[32,8] = 0.0
$arr[32,9] = 0.0
$arr[32,10] = 0.0
$arr[32,11] = 0.0
$arr[32,12] = 0.0
$arr[32,13] = 0.0
$arr[33,0] = 45740.34375
$arr[33,1] = 0.0
$arr[33,2] = 0.0
$arr[33,3] = 0.0
$arr[33,4] = 0.0
$arr[33,5] = 0.0
$arr[33,6] = 0.0
$arr[33,7] = 0.0
$arr[33,8] = 0.0
$arr[33,9] = 0.0
$arr[33,10] = 0.0
$arr[33,11] = 0.0
$arr[33,12] = 0.0
$arr[33,13] = 0.0
$arr[34,0] = 45740.35416666666
$arr[34,1] = 0.0
$arr[34,2] = 0.0
$arr[34,3] = 0.0
$arr[34,4] = 0.0
$arr[34,5] = 0.0
$arr[34,6] = 0.0
$arr[34,7] = 0.0
$arr[34,8] = 0.0
$arr[34,9] = 0.0
$arr[34,10] = 0.0
$arr[34,11] = 0.0
$arr[34,12] = 0.0
$arr[34,13] = 0.0
$arr[35,0] = 45740.36458333334
$arr[35,1] = 0.0
$arr[35,2] = 0.0
$arr[35,3] = 0.0
$arr[35,4] = 0.0
$arr[35,5] = 0.0
$arr[35,6] = 0.0
$arr[35,7] = 0.0
$arr[35,8] = 0.0
$arr[35,9] = 0.0
$arr[35,10] = 0.0
$arr[35,11] = 0.0
$arr[35,12] = 0.0
$arr[35,13] = 0.0
$arr[36,0] = 45740.375
$arr[36,1] = 0.0
$arr[36,2] = 0.0
$arr[36,3] = 0.0
$arr[36,4] = 0.0
$arr[36,5] = 0.0
$arr[36,6] = 0.0
$arr[36,7] = 0.0
$arr[36,8] = 0.0
$arr[36,9] = 0.0
$arr[36,10] = 0.0
$arr[36,11] = 0.0
$arr[36,12] = 0.0
$arr[36,13] = 0.0
$arr[37,0] = 45740.38541666666
$arr[37,1] = 0.0
$arr[37,2] = 0.0
$arr[37,3] = 0.0
$arr[37,4] = 0.0
$arr[37,5] = 0.0
$arr[37,6] = 0.0
$arr[37,7] = 0.0
$arr[37,8] = 0.0
$arr[37,9] = 0.0
$arr[37,10] = 0.0
$arr[37,11] = 0.0
$arr[37,12] = 0.0
$arr[37,13] = 0.0
$arr[38,0] = 45740.39583333334
$arr[38,1] = 0.0
$arr[38,2] = 0.0
$arr[38,3] = 0.0
$arr[38,4] = 0.0
$arr[38,5] = 0.0
$arr[38,6] = 0.0
$arr[38,7] = 0.0
$arr[38,8] = 0.0
$arr[38,9] = 0.0
$arr[38,10] = 0.0
$arr[38,11] = 0.0
$arr[38,12] = 0.0
$arr[38,13] = 0.0
$arr[39,0] = 45740.40625
$arr[39,1] = 0.0
$arr[39,2] = 0.0
$arr[39,3] = 0.0
$arr[39,4] = 0.0
$arr[39,5] = 0.0
$arr[39,6] = 0.0
$arr[39,7] = 0.0
$arr[39,8] = 0.0
$arr[39,9] = 0.0
$arr[39,10] = 0.0
$arr[39,11] = 0.0
$arr[39,12] = 0.0
$arr[39,13] = 0.0
$arr[40,0] = 45740.41666666666
$arr[40,1] = 0.0
$arr[40,2] = 0.0
$arr[40,3] = 0.0
$arr[40,4] = 0.0
$arr[40,5] = 0.0
$arr[40,6] = 0.0
$arr[40,7] = 0.0
$arr[40,8] = 0.0
$arr[40,9] = 0.0
$arr[40,10] = 0.0
$arr[40,11] = 0.0
$arr[40,12] = 0.0
$arr[40,13] = 0.0
$arr[41,0] = 45740.42708333334
$arr[41,1] = 0.0
$arr[41,2] = 0.0
$arr[41,3] = 0.0
$arr[41,4] = 0.0
$arr[41,5] = 0.0
$arr[41,6] = 0.0
$arr[41,7] = 0.0
$arr[41,8] = 0.0
$arr[41,9] = 0.0
$arr[41,10] = 0.0
$arr[41,11] = 0.0
$arr[41,12] = 0.0
$arr[41,13] = 0.0
$arr[42,0] = 45740.4375
$arr[42,1] = 0.0
$arr[42,2] = 0.0
$arr[42,3] = 0.0
$arr[42,4] = 0.0
$arr[42,5] = 0.0
$arr[42,6] = 0.0
$arr[42,7] = 0.0
$arr[42,8] = 0.0
$arr[42,9] = 0.0
$arr[42,10] = 0.0
$arr[42,11] = 0.0
$arr[42,12] = 0.0
$arr[42,13] = 0.0
$arr[43,0] = 45740.44791666666
$arr[43,1] = 0.0
$arr[43,2] = 0.0
$arr[43,3] = 0.0
$arr[43,4] = 0.0
$arr[43,5] = 0.0
$arr[43,6] = 0.0
$arr[43,7] = 0.0
$arr[43,8] = 0.0
$arr[43,9] = 0.0
$arr[43,10] = 0.0
$arr[43,11] = 0.0
$arr[43,12] = 0.0
$arr[43,13] = 0.0
$arr[44,0] = 45740.45833333334
$arr[44,1] = 0.0
$arr[44,2] = 0.0
$arr[44,3] = 0.0
$arr[44,4] = 0.0
$arr[44,5] = 0.0
$arr[44,6] = 0.0
$arr[44,7] = 0.0
$arr[44,8] = 0.0
$arr[44,9] = 0.0
$arr[44,10] = 0.0
$arr[44,11] = 0.0
$arr[44,12] = 0.0
$arr[44,13] = 0.0
$arr[45,0] = 45740.46875
$arr[45,1] = 0.0
$arr[45,2] = 0.0
$arr[45,3] = 0.0
$arr[45,4] = 0.0
$arr[45,5] = 0.0
$arr[45,6] = 0.0
$arr[45,7] = 0.0
$arr[45,8] = 0.0
$arr[45,9] = 0.0
$arr[45,10] = 0.0
$arr[45,11] = 0.0
$arr[45,12] = 0.0
$arr[45,13] = 0.0
$arr[46,0] = 45740.47916666666
$arr[46,1] = 0.0
$arr[46,2] = 0.0
$arr[46,3] = 0.0
$arr[46,4] = 0.0
$arr[46,5] = 0.0
$arr[46,6] = 0.0
$arr[46,7] = 0.0
$arr[46,8] = 0.0
$arr[46,9] = 0.0
$arr[46,10] = 0.0
$arr[46,11] = 0.0
$arr[46,12] = 0.0
$arr[46,13] = 0.0
$arr[47,0] = 45740.48958333334
$arr[47,1] = 0.0
$arr[47,2] = 0.0
$arr[47,3] = 0.0
$arr[47,4] = 0.0
$arr[47,5] = 0.0
$arr[47,6] = 0.0
$arr[47,7] = 0.0
$arr[47,8] = 0.0
$arr[47,9] = 0.0
$arr[47,10] = 0.0
$arr[47,11] = 0.0
$arr[47,12] = 0.0
$arr[47,13] = 0.0
$arr[48,0] = 45740.5
$arr[48,1] = 0.0
$arr[48,2] = 0.0
$arr[48,3] = 0.0
$arr[48,4] = 0.0
$arr[48,5] = 0.0
$arr[48,6] = 0.0
$arr[48,7] = 0.0
$arr[48,8] = 0.0
$arr[48,9] = 0.0
$arr[48,10] = 0.0
$arr[48,11] = 0.0
$arr[48,12] = 0.0
$arr[48,13] = 0.0
$arr[49,0] = 45740.51041666666
$arr[49,1] = 0.0
$arr[49,2] = 0.0
$arr[49,3] = 0.0
$arr[49,4] = 0.0
$arr[49,5] = 0.0
$arr[49,6] = 0.0
$arr[49,7] = 0.0
$arr[49,8] = 0.0
$arr[49,9] = 0.0
$arr[49,10] = 0.0
$arr[49,11] = 0.0
$arr[49,12] = 0.0
$arr[49,13] = 0.0
$arr[50,0] = 45740.52083333334
$arr[50,1] = 0.0
$arr[50,2] = 0.0
$arr[50,3] = 0.0
$arr[50,4] = 0.0
$arr[50,5] = 0.0
$arr[50,6] = 0.0
$arr[50,7] = 0.0
$arr[50,8] = 0.0
$arr[50,9] = 0.0
$arr[50,10] = 0.0
$arr[50,11] = 0.0
$arr[50,12] = 0.0
$arr[50,13] = 0.0
$arr[51,0] = 45740.53125
$arr[51,1] = 0.0
$arr[51,2] = 0.0
$arr[51,3] = 0.0
$arr[51,4] = 0.0
$arr[51,5] = 0.0
$arr[51,6] = 0.0
$arr[51,7] = 0.0
$arr[51,8] = 0.0
$arr[51,9] = 0.0
$arr[51,10] = 0.0
$arr[51,11] = 0.0
$arr[51,12] = 0.0
$arr[51,13] = 0.0
$arr[52,0] = 45740.54166666666
$arr[52,1] = 0.0
$arr[52,2] = 0.0
$arr[52,3] = 0.0
$arr[52,4] = 0.0
$arr[52,5] = 0.0
$arr[52,6] = 0.0
$arr[52,7] = 0.0
$arr[52,8] = 0.0
$arr[52,9] = 0.0
$arr[52,10] = 0.0
$arr[52,11] = 0.0
$arr[52,12] = 0.0
$arr[52,13] = 0.0
$arr[53,0] = 45740.55208333334
$arr[53,1] = 0.0
$arr[53,2] = 0.0
$arr[53,3] = 0.0
$arr[53,4] = 0.0
$arr[53,5] = 0.0
$arr[53,6] = 0.0
$arr[53,7] = 0.0
$arr[53,8] = 0.0
$arr[53,9] = 0.0
$arr[53,10] = 0.0
$arr[53,11] = 0.0
$arr[53,12] = 0.0
$arr[53,13] = 0.0
$arr[54,0] = 45740.5625
$arr[54,1] = 0.0
$arr[54,2] = 0.0
$arr[54,3] = 0.0
$arr[54,4] = 0.0
$arr[54,5] = 0.0
$arr[54,6] = 0.0
$arr[54,7] = 0.0
$arr[54,8] = 0.0
$arr[54,9] = 0.0
$arr[54,10] = 0.0
$arr[54,11] = 0.0
$arr[54,12] = 0.0
$arr[54,13] = 0.0
$arr[55,0] = 45740.57291666666
$arr[55,1] = 0.0
$arr[55,2] = 0.0
$arr[55,3] = 0.0
$arr[55,4] = 0.0
$arr[55,5] = 0.0
$arr[55,6] = 0.0
$arr[55,7] = 0.0
$arr[55,8] = 0.0
$arr[55,9] = 0.0
$arr[55,10] = 0.0
$arr[55,11] = 0.0
$arr[55,12] = 0.0
$arr[55,13] = 0.0
$arr[56,0] = 45740.58333333334
$arr[56,1] = 0.0
$arr[56,2] = 0.0
$arr[56,3] = 0.0
$arr[56,4] = 0.0
$arr[56,5] = 0.0
$arr[56,6] = 0.0
$arr[56,7] = 0.0
$arr[56,8] = 0.0
$arr[56,9] = 0.0
$arr[56,10] = 0.0
$arr[56,11] = 0.0
$arr[56,12] = 0.0
$arr[56,13] = 0.0
$arr[57,0] = 45740.59375
$arr[57,1] = 0.0
$arr[57,2] = 0.0
$arr[57,3] = 0.0
$arr[57,4] = 0.0
$arr[57,5] = 0.0
$arr[57,6] = 0.0
$arr[57,7] = 0.0
$arr[57,8] = 0.0
$arr[57,9] = 0.0
$arr[57,10] = 0.0
$arr[57,11] = 0.0
$arr[57,12] = 0.0
$arr[57,13] = 0.0
$arr[58,0] = 45740.60416666666
$arr[58,1] = 0.0
$arr[58,2] = 0.0
$arr[58,3] = 0.0
$arr[58,4] = 0.0
$arr[58,5] = 0.0
$arr[58,6] = 0.0
$arr[58,7] = 0.0
$arr[58,8] = 0.0
$arr[58,9] = 0.0
$arr[58,10] = 0.0
$arr[58,11] = 0.0
$arr[58,12] = 0.0
$arr[58,13] = 0.0
$arr[59,0] = 45740.61458333334
$arr[59,1] = 0.0
$arr[59,2] = 0.0
$arr[59,3] = 0.0
$arr[59,4] = 0.0
$arr[59,5] = 0.0
$arr[59,6] = 0.0
$arr[59,7] = 0.0
$arr[59,8] = 0.0
$arr[59,9] = 0.0
$arr[59,10] = 0.0
$arr[59,11] = 0.0
$arr[59,12] = 0.0
$arr[59,13] = 0.0
$arr[60,0] = 45740.625
$arr[60,1] = 0.0
$arr[60,2] = 0.0
$arr[60,3] = 0.0
$arr[60,4] = 0.0
$arr[60,5] = 0.0
$arr[60,6] = 0.0
$arr[60,7] = 0.0
$arr[60,8] = 0.0
$arr[60,9] = 0.0
$arr[60,10] = 0.0
$arr[60,11] = 0.0
$arr[60,12] = 0.0
$arr[60,13] = 0.0
$arr[61,0] = 45740.63541666666
$arr[61,1] = 0.0
$arr[61,2] = 0.0
$arr[61,3] = 0.0
$arr[61,4] = 0.0
$arr[61,5] = 0.0
$arr[61,6] = 0.0
$arr[61,7] = 0.0
$arr[61,8] = 0.0
$arr[61,9] = 0.0
$arr[61,10] = 0.0
$arr[61,11] = 0.0
$arr[61,12] = 0.0
$arr[61,13] = 0.0
$arr[62,0] = 45740.64583333334
$arr[62,1] = 0.0
$arr[62,2] = 0.0
$arr[62,3] = 0.0
$arr[62,4] = 0.0
$arr[62,5] = 0.0
$arr[62,6] = 0.0
$arr[62,7] = 0.0
$arr[62,8] = 0.0
$arr[62,9] = 0.0
$arr[62,10] = 0.0
$arr[62,11] = 0.0
$arr[62,12] = 0.0
$arr[62,13] = 0.0
$arr[63,0] = 45740.65625
$arr[63,1] = 0.0
$arr[63,2] = 0.0
$arr[63,3] = 0.0
$arr[63,4] = 0.0
$arr[63,5] = 0.0
$arr[63,6] = 0.0
$arr[63,7] = 0.0
$arr[63,8] = 0.0
$arr[63,9] = 0.0
$arr[63,10] = 0.0
$arr[63,11] = 0.0
$arr[63,12] = 0.0
$arr[63,13] = 0.0
$arr[64,0] = 45740.66666666666
$arr[64,1] = 0.0
$arr[64,2] = 0.0
$arr[64,3] = 0.0
$arr[64,4] = 0.0
$arr[64,5] = 0.0
$arr[64,6] = 0.0
$arr[64,7] = 0.0
$arr[64,8] = 0.0
$arr[64,9] = 0.0
$arr[64,10] = 0.0
$arr[64,11] = 0.0
$arr[64,12] = 0.0
$arr[64,13] = 0.0
$arr[65,0] = 45740.67708333334
$arr[65,1] = 0.0
$arr[65,2] = 0.0
$arr[65,3] = 0.0
$arr[65,4] = 0.0
$arr[65,5] = 0.0
$arr[65,6] = 0.0
$arr[65,7] = 0.0
$arr[65,8] = 0.0
$arr[65,9] = 0.0
$arr[65,10] = 0.0
$arr[65,11] = 0.0
$arr[65,12] = 0.0
$arr[65,13] = 0.0
$arr[66,0] = 45740.6875
$arr[66,1] = 0.0
$arr[66,2] = 0.0
$arr[66,3] = 0.0
$arr[66,4] = 0.0
$arr[66,5] = 0.0
$arr[66,6] = 0.0
$arr[66,7] = 0.0
$arr[66,8] = 0.0
$arr[66,9] = 0.0
$arr[66,10] = 0.0
$arr[66,11] = 0.0
$arr[66,12] = 0.0
$arr[66,13] = 0.0
$arr[67,0] = 45740.69791666666
$arr[67,1] = 0.0
$arr[67,2] = 0.0
$arr[67,3] = 0.0
$arr[67,4] = 0.0
$arr[67,5] = 0.0
$arr[67,6] = 0.0
$arr[67,7] = 0.0
$arr[67,8] = 0.0
$arr[67,9] = 0.0
$arr[67,10] = 0.0
$arr[67,11] = 0.0
$arr[67,12] = 0.0
$arr[67,13] = 0.0
$arr[68,0] = 45740.70833333334
$arr[68,1] = 0.0
$arr[68,2] = 0.0
$arr[68,3] = 0.0
$arr[68,4] = 0.0
$arr[68,5] = 0.0
$arr[68,6] = 0.0
$arr[68,7] = 0.0
$arr[68,8] = 0.0
$arr[68,9] = 0.0
$arr[68,10] = 0.0
$arr[68,11] = 0.0
$arr[68,12] = 0.0
$arr[68,13] = 0.0
$arr[69,0] = 45740.71875
$arr[69,1] = 0.0
$arr[69,2] = 0.0
$arr[69,3] = 0.0
$arr[69,4] = 0.0
$arr[69,5] = 0.0
$arr[69,6] = 0.0
$arr[69,7] = 0.0
$arr[69,8] = 0.0
$arr[69,9] = 0.0
$arr[69,10] = 0.0
$arr[69,11] = 0.0
$arr[69,12] = 0.0
$arr[69,13] = 0.0
$arr[70,0] = 45740.72916666666
$arr[70,1] = 0.0
$arr[70,2] = 0.0
$arr[70,3] = 0.0
$arr[70,4] = 0.0
$arr[70,5] = 0.0
$arr[70,6] = 0.0
$arr[70,7] = 0.0
$arr[70,8] = 0.0
$arr[70,9] = 0.0
$arr[70,10] = 0.0
$arr[70,11] = 0.0
$arr[70,12] = 0.0
$arr[70,13] = 0.0
$arr[71,0] = 45740.73958333334
$arr[71,1] = 0.0
$arr[71,2] = 0.0
$arr[71,3] = 0.0
$arr[71,4] = 0.0
$arr[71,5] = 0.0
$arr[71,6] = 0.0
$arr[71,7] = 0.0
$arr[71,8] = 0.0
$arr[71,9] = 0.0
$arr[71,10] = 0.0
$arr[71,11] = 0.0
$arr[71,12] = 0.0
$arr[71,13] = 0.0
$arr[72,0] = 45740.75
$arr[72,1] = 0.0
$arr[72,2] = 0.0
$arr[72,3] = 0.0
$arr[72,4] = 0.0
$arr[72,5] = 0.0
$arr[72,6] = 0.0
$arr[72,7] = 0.0
$arr[72,8] = 0.0
$arr[72,9] = 0.0
$arr[72,10] = 0.0
$arr[72,11] = 0.0
$arr[72,12] = 0.0
$arr[72,13] = 0.0
$arr[73,0] = 45740.76041666666
$arr[73,1] = 0.0
$arr[73,2] = 0.0
$arr[73,3] = 0.0
$arr[73,4] = 0.0
$arr[73,5] = 0.0
$arr[73,6] = 0.0
$arr[73,7] = 0.0
$arr[73,8] = 0.0
$arr[73,9] = 0.0
$arr[73,10] = 0.0
$arr[73,11] = 0.0
$arr[73,12] = 0.0
$arr[73,13] = 0.0
$arr[74,0] = 45740.77083333334
$arr[74,1] = 0.0
$arr[74,2] = 0.0
$arr[74,3] = 0.0
$arr[74,4] = 0.0
$arr[74,5] = 0.0
$arr[74,6] = 0.0
$arr[74,7] = 0.0
$arr[74,8] = 0.0
$arr[74,9] = 0.0
$arr[74,10] = 0.0
$arr[74,11] = 0.0
$arr[74,12] = 0.0
$arr[74,13] = 0.0
$arr[75,0] = 45740.78125
$arr[75,1] = 0.0
$arr[75,2] = 0.0
$arr[75,3] = 0.0
$arr[75,4] = 0.0
$arr[75,5] = 0.0
$arr[75,6] = 0.0
$arr[75,7] = 0.0
$arr[75,8] = 0.0
$arr[75,9] = 0.0
$arr[75,10] = 0.0
$arr[75,11] = 0.0
$arr[75,12] = 0.0
$arr[75,13] = 0.0
$arr[76,0] = 45740.79166666666
$arr[76,1] = 0.0
$arr[76,2] = 0.0
$arr[76,3] = 0.0
$arr[76,4] = 0.0
$arr[76,5] = 0.0
$arr[76,6] = 0.0
$arr[76,7] = 0.0
$arr[76,8] = 0.0
$arr[76,9] = 0.0
$arr[76,10] = 0.0
$arr[76,11] = 0.0
$arr[76,12] = 0.0
$arr[76,13] = 0.0
$arr[77,0] = 45740.80208333334
$arr[77,1] = 0.0
$arr[77,2] = 0.0
$arr[77,3] = 0.0
$arr[77,4] = 0.0
$arr[77,5] = 0.0
$arr[77,6] = 0.0
$arr[77,7] = 0.0
$arr[77,8] = 0.0
$arr[77,9] = 0.0
$arr[77,10] = 0.0
$arr[77,11] = 0.0
$arr[77,12] = 0.0
$arr[77,13] = 0.0
$arr[78,0] = 45740.8125
$arr[78,1] = 0.0
$arr[78,2] = 0.0
$arr[78,3] = 0.0
$arr[78,4] = 0.0
$arr[78,5] = 0.0
$arr[78,6] = 0.0
$arr[78,7] = 0.0
$arr[78,8] = 0.0
$arr[78,9] = 0.0
$arr[78,10] = 0.0
$arr[78,11] = 0.0
$arr[78,12] = 0.0
$arr[78,13] = 0.0
$arr[79,0] = 45740.82291666666
$arr[79,1] = 0.0
$arr[79,2] = 0.0
$arr[79,3] = 0.0
$arr[79,4] = 0.0
$arr[79,5] = 0.0
$arr[79,6] = 0.0
$arr[79,7] = 0.0
$arr[79,8] = 0.0
$arr[79,9] = 0.0
$arr[79,10] = 0.0
$arr[79,11] = 0.0
$arr[79,12] = 0.0
$arr[79,13] = 0.0
$arr[80,0] = 45740.83333333334
$arr[80,1] = 0.0
$arr[80,2] = 0.0
$arr[80,3] = 0.0
$arr[80,4] = 0.0
$arr[80,5] = 0.0
$arr[80,6] = 0.0
$arr[80,7] = 0.0
$arr[80,8] = 0.0
$arr[80,9] = 0.0
$arr[80,10] = 0.0
$arr[80,11] = 0.0
$arr[80,12] = 0.0
$arr[80,13] = 0.0
$arr[81,0] = 45740.84375
$arr[81,1] = 0.0
$arr[81,2] = 0.0
$arr[81,3] = 0.0
$arr[81,4] = 0.0
$arr[81,5] = 0.0
$arr[81,6] = 0.0
$arr[81,7] = 0.0
$arr[81,8] = 0.0
$arr[81,9] = 0.0
$arr[81,10] = 0.0
$arr[81,11] = 0.0
$arr[81,12] = 0.0
$arr[81,13] = 0.0
$arr[82,0] = 45740.85416666666
$arr[82,1] = 0.0
$arr[82,2] = 0.0
$arr[82,3] = 0.0
$arr[82,4] = 0.0
$arr[82,5] = 0.0
$arr[82,6] = 0.0
$arr[82,7] = 0.0
$arr[82,8] = 0.0
$arr[82,9] = 0.0
$arr[82,10] = 0.0
$arr[82,11] = 0.0
$arr[82,12] = 0.0
$arr[82,13] = 0.0
$arr[83,0] = 45740.86458333334
$arr[83,1] = 0.0
$arr[83,2] = 0.0
$arr[83,3] = 0.0
$arr[83,4] = 0.0
$arr[83,5] = 0.0
$arr[83,6] = 0.0
$arr[83,7] = 0.0
$arr[83,8] = 0.0
$arr[83,9] = 0.0
$arr[83,10] = 0.0
$arr[83,11] = 0.0
$arr[83,12] = 0.0
$arr[83,13] = 0.0
$arr[84,0] = 45740.875
$arr[84,1] = 0.0
$arr[84,2] = 0.0
$arr[84,3] = 0.0
$arr[84,4] = 0.0
$arr[84,5] = 0.0
$arr[84,6] = 0.0
$arr[84,7] = 0.0
$arr[84,8] = 0.0
$arr[84,9] = 0.0
$arr[84,10] = 0.0
$arr[84,11] = 0.0
$arr[84,12] = 0.0
$arr[84,13] = 0.0
$arr[85,0] = 45740.88541666666
$arr[85,1] = 0.0
$arr[85,2] = 0.0
$arr[85,3] = 0.0
$arr[85,4] = 0.0
$arr[85,5] = 0.0
$arr[85,6] = 0.0
$arr[85,7] = 0.0
$arr[85,8] = 0.0
$arr[85,9] = 0.0
$arr[85,10] = 0.0
$arr[85,11] = 0.0
$arr[85,12] = 0.0
$arr[85,13] = 0.0
$arr[86,0] = 45740.89583333334
$arr[86,1] = 0.0
$arr[86,2] = 0.0
$arr[86,3] = 0.0
$arr[86,4] = 0.0
$arr[86,5] = 0.0
$arr[86,6] = 0.0
$arr[86,7] = 0.0
$arr[86,8] = 0.0
$arr[86,9] = 0.0
$arr[86,10] = 0.0
$arr[86,11] = 0.0
$arr[86,12] = 0.0
$arr[86,13] = 0.0
$arr[87,0] = 45740.90625
$arr[87,1] = 0.0
$arr[87,2] = 0.0
$arr[87,3] = 0.0
$arr[87,4] = 0.0
$arr[87,5] = 0.0
$arr[87,6] = 0.0
$arr[87,7] = 0.0
$arr[87,8] = 0.0
$arr[87,9] = 0.0
$arr[87,10] = 0.0
$arr[87,11] = 0.0
$arr[87,12] = 0.0
$arr[87,13] = 0.0
$arr[88,0] = 45740.91666666666
$arr[88,1] = 0.0
$arr[88,2] = 0.0
$arr[88,3] = 0.0
$arr[88,4] = 0.0
$arr[88,5] = 0.0
$arr[88,6] = 0.0
$arr[88,7] = 0.0
$arr[88,8] = 0.0
$arr[88,9] = 0.0
$arr[88,10] = 0.0
$arr[88,11] = 0.0
$arr[88,12] = 0.0
$arr[88,13] = 0.0
$arr[89,0] = 45740.92708333334
$arr[89,1] = 0.0
$arr[89,2] = 0.0
$arr[89,3] = 0.0
$arr[89,4] = 0.0
$arr[89,5] = 0.0
$arr[89,6] = 0.0
$arr[89,7] = 0.0
$arr[89,8] = 0.0
$arr[89,9] = 0.0
$arr[89,10] = 0.0
$arr[89,11] = 0.0
$arr[89,12] = 0.0
$arr[89,13] = 0.0
$arr[90,0] = 45740.9375
$arr[90,1] = 0.0
$arr[90,2] = 0.0
$arr[90,3] = 0.0
$arr[90,4] = 0.0
$arr[90,5] = 0.0
$arr[90,6] = 0.0
$arr[90,7] = 0.0
$arr[90,8] = 0.0
$arr[90,9] = 0.0
$arr[90,10] = 0.0
$arr[90,11] = 0.0
$arr[90,12] = 0.0
$arr[90,13] = 0.0
$arr[91,0] = 45740.94791666666
$arr[91,1] = 0.0
$arr[91,2] = 0.0
$arr[91,3] = 0.0
$arr[91,4] = 0.0
$arr[91,5] = 0.0
$arr[91,6] = 0.0
$arr[91,7] = 0.0
$arr[91,8] = 0.0
$arr[91,9] = 0.0
$arr[91,10] = 0.0
$arr[91,11] = 0.0
$arr[91,12] = 0.0
$arr[91,13] = 0.0
$arr[92,0] = 45740.95833333334
$arr[92,1] = 0.0
$arr[92,2] = 0.0
$arr[92,3] = 0.0
$arr[92,4] = 0.0
$arr[92,5] = 0.0
$arr[92,6] = 0.0
$arr[92,7] = 0.0
$arr[92,8] = 0.0
$arr[92,9] = 0.0
$arr[92,10] = 0.0
$arr[92,11] = 0.0
$arr[92,12] = 0.0
$arr[92,13] = 0.0
$arr[93,0] = 45740.96875
$arr[93,1] = 0.0
$arr[93,2] = 0.0
$arr[93,3] = 0.0
$arr[93,4] = 0.0
$arr[93,5] = 0.0
$arr[93,6] = 0.0
$arr[93,7] = 0.0
$arr[93,8] = 0.0
$arr[93,9] = 0.0
$arr[93,10] = 0.0
$arr[93,11] = 0.0
$arr[93,12] = 0.0
$arr[93,13] = 0.0
$arr[94,0] = 45740.97916666666
$arr[94,1] = 0.0
$arr[94,2] = 0.0
$arr[94,3] = 0.0
$arr[94,4] = 0.0
$arr[94,5] = 0.0
$arr[94,6] = 0.0
$arr[94,7] = 0.0
$arr[94,8] = 0.0
$arr[94,9] = 0.0
$arr[94,10] = 0.0
$arr[94,11] = 0.0
$arr[94,12] = 0.0
$arr[94,13] = 0.0
$arr[95,0] = 45740.98958333334
$arr[95,1] = 0.0
$arr[95,2] = 0.0
$arr[95,3] = 0.0
$arr[95,4] = 0.0
$arr[95,5] = 0.0
$arr[95,6] = 0.0
$arr[95,7] = 0.0
$arr[95,8] = 0.0
$arr[95,9] = 0.0
$arr[95,10] = 0.0
$arr[95,11] = 0.0
$arr[95,12] = 0.0
$arr[95,13] = 0.0
$ws.Range("A2:N97").Value = $arr
